$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("ALC")
$ws2 = $wb.Worksheets.Item("ARM")
$ws3 = $wb.Worksheets.Item("CRP")
$ws4 = $wb.Worksheets.Item("CUL")
$ws5 = $wb.Worksheets.Item("GSM")
$ws6 = $wb.Worksheets.Item("LTW")
$ws7 = $wb.Worksheets.Item("WVR")

# --- ALC ---
$ws1.Range("H6").Value = 1077.2
$ws1.Range("I6").Value = 1077.2
$ws1.Range("K6").Value = 3231.6
$ws1.Range("M6").Value = -3119.6
$ws1.Range("H15").Value = 1961.5128
$ws1.Range("I15").Value = 1961.5128
$ws1.Range("K15").Value = 5884.538399999999
$ws1.Range("M15").Value = -5715.538399999999
$ws1.Range("H20").Value = 0
$ws1.Range("J20").Value = 0
$ws1.Range("L20").Value = 0
$ws1.Range("N20").ClearContents()
$ws1.Range("H33").Value = 262.41934
$ws1.Range("I33").Value = 268.03333
$ws1.Range("K33").Value = 268.03333
$ws1.Range("M33").Value = -39.03332999999998
$ws1.Range("H35").Value = 0
$ws1.Range("J35").Value = 0
$ws1.Range("L35").Value = 0
$ws1.Range("N35").ClearContents()
$ws1.Range("H64").Value = 3996.6667
$ws1.Range("I64").Value = 3994
$ws1.Range("K64").Value = 3994
$ws1.Range("M64").Value = -3746
$ws1.Range("H67").Value = 3996.6667
$ws1.Range("I67").Value = 3994
$ws1.Range("K67").Value = 3994
$ws1.Range("M67").Value = -3136
$ws1.Range("H74").Value = 15158.529
$ws1.Range("I74").Value = 15913.134
$ws1.Range("J74").Value = 9499
$ws1.Range("K74").Value = 15913.134
$ws1.Range("L74").Value = 9499
$ws1.Range("M74").Value = -14977.134
$ws1.Range("N74").Value = -11371
$ws1.Range("H77").Value = 15158.529
$ws1.Range("I77").Value = 15913.134
$ws1.Range("J77").Value = 9499
$ws1.Range("K77").Value = 79565.67
$ws1.Range("L77").Value = 47495
$ws1.Range("M77").Value = -74885.67
$ws1.Range("N77").Value = -56855
$ws1.Range("H111").Value = 15125
$ws1.Range("I111").Value = 15125
$ws1.Range("J111").Value = 0
$ws1.Range("K111").Value = 45375
$ws1.Range("L111").Value = 0
$ws1.Range("M111").Value = -42308
$ws1.Range("N111").ClearContents()
$ws1.Range("H132").Value = 4205.3716
$ws1.Range("I132").Value = 4311.4116
$ws1.Range("K132").Value = 12934.2348
$ws1.Range("M132").Value = -10404.2348

# --- ARM ---
$ws2.Range("H31").Value = 7818.6
$ws2.Range("I31").Value = 7818.6
$ws2.Range("K31").Value = 7818.6
$ws2.Range("M31").Value = -7524.6
$ws2.Range("H32").Value = 1176713.6
$ws2.Range("I32").Value = 562080.9
$ws2.Range("K32").Value = 562080.9
$ws2.Range("M32").Value = -561793.9
$ws2.Range("H74").Value = 4598.1704
$ws2.Range("I74").Value = 2302.7932
$ws2.Range("J74").Value = 8296.277
$ws2.Range("K74").Value = 2302.7932
$ws2.Range("L74").Value = 8296.277
$ws2.Range("M74").Value = -1428.7932
$ws2.Range("N74").Value = -10044.277
$ws2.Range("H77").Value = 4598.1704
$ws2.Range("I77").Value = 2302.7932
$ws2.Range("J77").Value = 8296.277
$ws2.Range("K77").Value = 11513.966
$ws2.Range("L77").Value = 41481.385
$ws2.Range("M77").Value = -7145.966
$ws2.Range("N77").Value = -50217.385

# --- CRP ---
$ws3.Range("H31").Value = 2909298.5
$ws3.Range("I31").Value = 1559.4375
$ws3.Range("J31").Value = 4632403
$ws3.Range("K31").Value = 1559.4375
$ws3.Range("L31").Value = 4632403
$ws3.Range("M31").Value = -1264.4375
$ws3.Range("N31").Value = -4632993
$ws3.Range("H34").Value = 2909298.5
$ws3.Range("I34").Value = 1559.4375
$ws3.Range("J34").Value = 4632403
$ws3.Range("K34").Value = 1559.4375
$ws3.Range("L34").Value = 4632403
$ws3.Range("M34").Value = -1357.4375
$ws3.Range("N34").Value = -4632807
$ws3.Range("H50").Value = 54795
$ws3.Range("J50").Value = 54795
$ws3.Range("L50").Value = 54795
$ws3.Range("N50").Value = -56045
$ws3.Range("H58").Value = 2600
$ws3.Range("I58").Value = 1733.3334
$ws3.Range("K58").Value = 1733.3334
$ws3.Range("M58").Value = -1530.3334
$ws3.Range("H86").Value = 5248.75
$ws3.Range("I86").Value = 2003.5
$ws3.Range("J86").Value = 8494
$ws3.Range("K86").Value = 2003.5
$ws3.Range("L86").Value = 8494
$ws3.Range("M86").Value = -880.5
$ws3.Range("N86").Value = -10740
$ws3.Range("H89").Value = 5248.75
$ws3.Range("I89").Value = 2003.5
$ws3.Range("J89").Value = 8494
$ws3.Range("K89").Value = 10017.5
$ws3.Range("L89").Value = 42470
$ws3.Range("M89").Value = -4401.5
$ws3.Range("N89").Value = -53702
$ws3.Range("H134").Value = 2502.5813
$ws3.Range("I134").Value = 2458.8157
$ws3.Range("J134").Value = 2835.2
$ws3.Range("K134").Value = 7376.4471
$ws3.Range("L134").Value = 8505.599999999999
$ws3.Range("M134").Value = -4841.4471
$ws3.Range("N134").Value = -13575.6
$ws3.Range("H136").Value = 2600
$ws3.Range("I136").Value = 1733.3334
$ws3.Range("K136").Value = 5200.0002
$ws3.Range("M136").Value = -2650.0002
$ws3.Range("H138").Value = 99000
$ws3.Range("J138").Value = 99000
$ws3.Range("L138").Value = 99000
$ws3.Range("N138").Value = -109280

# --- CUL ---
$ws4.Range("H5").Value = 617.36365
$ws4.Range("I5").Value = 392.26666
$ws4.Range("J5").Value = 1099.7142
$ws4.Range("K5").Value = 1176.79998
$ws4.Range("L5").Value = 3299.1426
$ws4.Range("M5").Value = -1064.79998
$ws4.Range("N5").Value = -3523.1426
$ws4.Range("H68").Value = 1284524.9
$ws4.Range("I68").Value = 1926.5
$ws4.Range("J68").Value = 1473142.2
$ws4.Range("K68").Value = 5779.5
$ws4.Range("L68").Value = 4419426.6
$ws4.Range("M68").Value = -4968.5
$ws4.Range("N68").Value = -4421048.6
$ws4.Range("H71").Value = 1284524.9
$ws4.Range("I71").Value = 1926.5
$ws4.Range("J71").Value = 1473142.2
$ws4.Range("K71").Value = 17338.5
$ws4.Range("L71").Value = 13258279.8
$ws4.Range("M71").Value = -13282.5
$ws4.Range("N71").Value = -13266391.8
$ws4.Range("H81").Value = 2624.75
$ws4.Range("J81").Value = 4500
$ws4.Range("L81").Value = 13500
$ws4.Range("N81").Value = -15746
$ws4.Range("H84").Value = 2624.75
$ws4.Range("J84").Value = 4500
$ws4.Range("L84").Value = 40500
$ws4.Range("N84").Value = -51732
$ws4.Range("H129").Value = 63169.445
$ws4.Range("J129").Value = 112999
$ws4.Range("L129").Value = 338997
$ws4.Range("N129").Value = -348997
$ws4.Range("H135").Value = 617.36365
$ws4.Range("I135").Value = 392.26666
$ws4.Range("J135").Value = 1099.7142
$ws4.Range("K135").Value = 3530.39994
$ws4.Range("L135").Value = 9897.427799999999
$ws4.Range("M135").Value = -995.3999400000002
$ws4.Range("N135").Value = -14967.4278

# --- GSM ---
$ws5.Range("H136").Value = 87019.266
$ws5.Range("J136").Value = 87019.266
$ws5.Range("L136").Value = 261057.798
$ws5.Range("N136").Value = -266157.798

# --- LTW ---
$ws6.Range("H16").Value = 1696.762
$ws6.Range("I16").Value = 1445.8125
$ws6.Range("K16").Value = 1445.8125
$ws6.Range("M16").Value = -1275.8125
$ws6.Range("H40").Value = 89535.64
$ws6.Range("I40").Value = 151749.88
$ws6.Range("K40").Value = 151749.88
$ws6.Range("M40").Value = -151613.88
$ws6.Range("H61").Value = 6505.3887
$ws6.Range("I61").Value = 6436.5
$ws6.Range("K61").Value = 6436.5
$ws6.Range("M61").Value = -6234.5
$ws6.Range("H93").Value = 2491.4736
$ws6.Range("J93").Value = 2345.182
$ws6.Range("L93").Value = 2345.182
$ws6.Range("N93").Value = -4841.182
$ws6.Range("H112").Value = 57386.23
$ws6.Range("J112").Value = 57386.23
$ws6.Range("L112").Value = 57386.23
$ws6.Range("N112").Value = -60340.23
$ws6.Range("H113").Value = 6505.3887
$ws6.Range("I113").Value = 6436.5
$ws6.Range("K113").Value = 6436.5
$ws6.Range("M113").Value = -4266.5

# --- WVR ---
$ws7.Range("H30").Value = 16673818
$ws7.Range("I30").Value = 6669.3335
$ws7.Range("K30").Value = 6669.3335
$ws7.Range("M30").Value = -6562.3335
$ws7.Range("H70").Value = 39177.8
$ws7.Range("J70").Value = 38996.332
$ws7.Range("L70").Value = 38996.332
$ws7.Range("N70").Value = -39626.332
$ws7.Range("H73").Value = 39177.8
$ws7.Range("J73").Value = 38996.332
$ws7.Range("L73").Value = 38996.332
$ws7.Range("N73").Value = -41180.332
$ws7.Range("H122").Value = 8623019
$ws7.Range("I122").Value = 2582.5715
$ws7.Range("J122").Value = 31251664
$ws7.Range("K122").Value = 7747.7145
$ws7.Range("L122").Value = 93754992
$ws7.Range("M122").Value = -5297.7145
$ws7.Range("N122").Value = -93759892
